{"js": "// Apply the 5 text replacements described by the diff.\n// Each pair is [oldText, newText]; oldText occurs exactly once in the\n// document body, so a plain case-sensitive search() is sufficient.\nconst replacements = [\n  [\n    \"Que el correo y contrase\u00f1as del usuario est\u00e9n en la base de datos\",\n    \"Que el rut y contrase\u00f1as del usuario est\u00e9n en la base de datos\"\n  ],\n  [\n    \"El sistema desplegar\u00e1 en m\u00f3vil la geolocalizaci\u00f3n y en web el cat\u00e1logo de productos.\",\n    \"El sistema desplegar\u00e1 el menu de inicio\"\n  ],\n  [\n    \"El Proveedor/administrador selecciona el icono de perfil.\",\n    \"El Proveedor/administrador selecciona perfil.\"\n  ],\n  [\n    \"El sistema despliega la interfaz de \u00bfquieres vender tus productos?\",\n    \"El sistema despliega la interfaz de registro\"\n  ],\n  [\n    \"El correo del usuario no es v\u00e1lido (no existe).\",\n    \"El rut no es v\u00e1lido (no existe).\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const result of results.items) {\n    result.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the 5 text replacements described by the diff using Find/Replace\n# against the whole document story. Each old string occurs exactly once,\n# so MatchCase + non-wildcard Find.Execute with Replace:=wdReplaceAll (2)\n# safely targets only the intended run.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"Que el correo y contrase\u00f1as del usuario est\u00e9n en la base de datos\", \"Que el rut y contrase\u00f1as del usuario est\u00e9n en la base de datos\"),\n    @(\"El sistema desplegar\u00e1 en m\u00f3vil la geolocalizaci\u00f3n y en web el cat\u00e1logo de productos.\", \"El sistema desplegar\u00e1 el menu de inicio\"),\n    @(\"El Proveedor/administrador selecciona el icono de perfil.\", \"El Proveedor/administrador selecciona perfil.\"),\n    @(\"El sistema despliega la interfaz de \u00bfquieres vender tus productos?\", \"El sistema despliega la interfaz de registro\"),\n    @(\"El correo del usuario no es v\u00e1lido (no existe).\", \"El rut no es v\u00e1lido (no existe).\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $rng = $d.Content\n    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $findText\"\n    }\n}\n"}
